$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.73%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'27.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.66%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'4.862"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.02%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.06394"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.15%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'6.956"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.00%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.184"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-8.79%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.8755"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.47%"
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'3.54%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.05074"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.43%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07577"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.43%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.02956"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.32%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.08982"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.61%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001567"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.21%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0006388"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.93%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.006181"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'3.48%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.475"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.62%"
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'-0.37%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'-1.37%"
$ws.Range("E19").Style = "Normal"

$ws.Range("E21").Value = "'1.83%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'3.907"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.08%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04412"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.18%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D25").Value = "'0.001177"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.16%"
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'-9.44%"
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'-0.04%"
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'14.64%"
$ws.Range("E28").Style = "Normal"

$ws.Range("D40").Value = "'0.04163"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.15%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.006870"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.77%"
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'0.61%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002189"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'4.74%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01182"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.94%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005202"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-2.00%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.679"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'15.47%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.01853"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-7.34%"
$ws.Range("E47").Style = "Normal"
